$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the letterhead rows (rows 1-3) - shifts everything up
$ws.Range("A1:A3").EntireRow.Delete()

# Delete the trailing footnote row (previously row 16, now row 13 after the shift above)
$ws.Range("A13:A13").EntireRow.Delete()
